# Update "想去人数" (number of people interested) figures that changed
# between scrape runs, as described by the commit "Update gh-pages to
# output generated at 456a3b4".
#
# Sheet "展览" (Exhibition) and sheet "全部类型" (All types) both contain
# the same underlying rows (the latter is a superset), so the same F-column
# values need to be bumped in both places.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5401
$ws1.Range("F6").Value = 815
$ws1.Range("F7").Value = 13
$ws1.Range("F8").Value = 333
$ws1.Range("F9").Value = 14

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5401
$ws4.Range("F6").Value = 815
$ws4.Range("F7").Value = 13
$ws4.Range("F9").Value = 333
$ws4.Range("F10").Value = 14
